$d = $word.ActiveDocument
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*debemos conocer*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $r = $target.Range
    $r.Find.Execute("que", $true, $true, $false, $false, $false, $true, 1, $false, "qué", 2)
}
